# Adds season record columns (Wins / Losses / Ties) to the roster sheet.
# New headers go in AD1:AF1 (matching the header style used by the rest of
# row 1), and every data row (2-44) gets the team's season record:
# 85 wins, 77 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold font, border,
# centered/top aligned) by copying the format from the existing last
# header cell (AC1) onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record for every player row
$ws.Range("AD2:AD44").Value = 85
$ws.Range("AE2:AE44").Value = 77
$ws.Range("AF2:AF44").Value = 0
